# Update, report! (Update slide, báo cáo!)
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder fields: "7/12/2015" -> "9/12/2015"
#    (slide master, every slide layout, and the notes master)
# ---------------------------------------------------------------------
$newDate = "9/12/2015"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# ---------------------------------------------------------------------
# 2) Slide 19 ("TextBox 5"): merge the two runs
#    "Phần giao diện " + "chọn các chứn năng quản lý."
#    into a single run.
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
for ($i = 1; $i -le $s19.Shapes.Count; $i++) {
    $shp = $s19.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 5") {
        $tr = $shp.TextFrame.TextRange
        # force a real text diff so the runtime merges the text into one run,
        # then set the final (unchanged) text.
        $tr.Text = "placeholder reset"
        $tr.Text = "Phần giao diện chọn các chứn năng quản lý."
    }
}

# ---------------------------------------------------------------------
# 3) Slide 24 ("Content Placeholder 2"):
#    - split "...xem sản phẩm." into "...xem sản phẩm" + "."
#    - add a new bullet paragraph after it:
#      "Phát triển ứng dụng thu ngân trên các nền tảng di động."
# ---------------------------------------------------------------------
$s24 = $p.Slides.Item(24)
for ($i = 1; $i -le $s24.Shapes.Count; $i++) {
    $shp = $s24.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder 2") {
        $tr = $shp.TextFrame.TextRange

        $para4 = $tr.Paragraphs(4, 1)
        $para4.Text = "Phát triển thêm chức năng để người mua cũng có thể xem sản phẩmZZZ"
        $para4b = $tr.Paragraphs(4, 1)
        $para4b.Text = "Phát triển thêm chức năng để người mua cũng có thể xem sản phẩm."

        $para4c = $tr.Paragraphs(4, 1)
        $para4c.InsertAfter("`rPhát triển ứng dụng thu ngân trên các nền tảng di động.")
    }
}
